# Zombono version scheduling sheet: bump to 0.0.11
# Inserts 7 new "Bugfix" tasks right after the existing two bugfix rows
# (rows 5-11), pushing the remaining task rows down; moves the trailing
# "NOT FINISHED - MORE SOON" marker row down to match; updates the
# dimension / active selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row content (A = Task, B = Area) for rows 3..25 -----------------
# Rows 1 and 2 (title + header) are unchanged.

$rows = @{
    3  = @('Fix: Warehouse ramps too steep', 'Bugfix')
    4  = @('Fix " velocity increasing but not speed" (prediction miss bug when hitting wall at specific angle sometimes)', 'Bugfix')
    5  = @('Fix game connecting to server during intro', 'Bugfix')
    6  = @('Fix relatve velocity being added multiple times', 'Bugfix')
    7  = @('func_train in subway fucked up', 'Bugfix')
    8  = @('Fucked up chair collision', 'Bugfix')
    9  = @('Zombification: you can get stuck', 'Bugfix')
    10 = @('Green pipe, add poster there', 'Bugfix')
    11 = @('Game does not clear entities on 2nd entry into a map', 'Bugfix')
    12 = @('Add team door brush entity', 'Feature')
    13 = @('Complete Release Generation Tool', 'Engineering')
    14 = @('Add kill feed', 'Feature')
    15 = @('z_tdm_warehouse - Second Floor', 'Feature')
    16 = @('Make Master servers work', 'Feature, Netservuces')
    17 = @('BrowseServersUI', 'Feature, Netservices')
    18 = @('Text Engine - scaled coordinates', 'Feature')
    19 = @('Start Waves mode programming', 'Feature')
    20 = @('Properly split out client.h, server.h', 'Refactoring')
    21 = @('Allow people to see what team a player is', 'Feature')
    22 = @('cl_console_line_length', 'Feature')
    23 = @('z_waves_port working (THEY COME FROM THE SEA/!?!?!?!?!) - z_tdm_spire finished', 'Content')
    24 = @('Finish z_warehouse easter egg', 'Content')
    25 = @('NOT FINISHED - MORE SOON', $null)
}

for ($r = 3; $r -le 25; $r++) {
    $pair = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    if ($pair[1] -ne $null) {
        $ws.Cells.Item($r, 2).Value = $pair[1]
    } else {
        $ws.Cells.Item($r, 2).Value = ""
    }
}

# --- "Completion date" placeholder cells in column D ------------------
# These keep their special number formats, now anchored at the same row
# numbers as before the insert (10, 12, 15, 16, 17, 20); copy the
# formatting straight from the cells that already carry it so we reuse
# the existing style entries instead of minting new ones.

$ws.Cells.Item(10, 4).Copy()
$ws.Cells.Item(10, 4).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(12, 4).Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(15, 4).Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4122) | Out-Null

# 16, 17 and 20 already carried the "short date" style before the edit,
# so nothing further is needed there - they keep their formatting as-is.

# Row 18 no longer holds the bold "NOT FINISHED" marker - restore it to
# plain (non-bold) text and default row height.
$ws.Cells.Item(18, 1).Font.Bold = $false
$ws.Rows(18).EntireRow.AutoFit()

# Row 25 is now the bold "NOT FINISHED - MORE SOON" marker row.
$ws.Cells.Item(25, 1).Font.Bold = $true
$ws.Rows(25).RowHeight = 15

$ws.Application.CutCopyMode = $false

# --- Sheet-level bookkeeping -------------------------------------------
$ws.Range("A1:D25").Select() | Out-Null
$ws.Range("C11").Select()
